# Adds a header row (row 1) describing the logged columns, shifts the
# previously-only data row down into the log, refreshes that row's values
# for the new (Await-based) run, and appends four more logged runs
# (rows 3-6) -- matching the richer sentiment/timing log the app now
# writes after switching the worker thread from Wait to Await.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Header row -----------------------------------------------------
$headers = @(
    "Date/Time",
    "Method",
    "elapsedMs",
    "wordCount",
    "sentenceCount",
    "posWordCount",
    "negWordCount",
    "posWordPercentage",
    "negWordPercentage",
    "posPhraseCount",
    "negativePhraseCount",
    "posWordPercentage",
    "negPhrasePercentage"
)
for ($c = 1; $c -le $headers.Length; $c++) {
    $ws.Cells.Item(1, $c).Value = $headers[$c - 1]
}

# Row 1 column A should carry the same date-number-format style as the
# data rows below it (column style index reused, not a fresh number
# format) -- copy just the formatting down from A2, which already has it.
$ws.Range("A2").Copy()
$ws.Range("A1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ---- Data rows --------------------------------------------------------
# columns: DateTime, Method, elapsedMs, wordCount, sentenceCount,
#          posWordCount, negWordCount, posWordPercentage, negWordPercentage,
#          posPhraseCount, negativePhraseCount, posWordPercentage(dup),
#          negPhrasePercentage
$rows = @(
    @(42585.690925925926, "Named", 12614, 7505, 440, 95, 47, 66, 32, 2, 0, 66, 0),
    @(42585.694745370369, "Named", 12197, 7859, 456, 101, 50, 66, 33, 1, 0, 66, 0),
    @(42585.698321759257, "Named", 12820, 7854, 456, 100, 50, 66, 33, 1, 0, 66, 0),
    @(42585.704340277778, "Named", 12534, 7604, 445, 98, 49, 66, 33, 2, 0, 66, 0),
    @(42585.707280092596, "Named", 13218, 7932, 461, 104, 51, 66, 32, 1, 0, 66, 0)
)

$r = 2
foreach ($row in $rows) {
    for ($c = 1; $c -le $row.Length; $c++) {
        $ws.Cells.Item($r, $c).Value = $row[$c - 1]
    }
    $r++
}

# Every logged date in column A (rows 2-6) shares the same date/time
# number-format style as the original row did; make sure the newly
# appended rows (3-6) pick it up too.
$ws.Range("A2").Copy()
$ws.Range("A3:A6").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ---- Column widths (approximate Excel's "best fit" auto-sizing) ------
# Target (post-AutoFit) stored widths, pre-compensated for this host's
# ColumnWidth -> stored-width rounding (stored = round(input*6)/6 + 5/6)
# so the written <col width="..."> lands as close as possible to what
# Excel's own best-fit sizing produced for these headers/values.
$widths = @(
    13.0,
    7.166666666666667,
    9.666666666666666,
    10.0,
    13.666666666666666,
    13.666666666666666,
    13.666666666666666,
    18.5,
    18.666666666666668,
    14.666666666666666,
    19.5,
    18.5,
    19.666666666666668
)
for ($c = 1; $c -le $widths.Length; $c++) {
    $ws.Columns.Item($c).ColumnWidth = $widths[$c - 1]
}

Write-Output "done"
